# Insert two new "Author" paragraphs ("Ben Jarman" and "Catherine Heard")
# right after the Subtitle paragraph and before the Date paragraph.

$d = $word.ActiveDocument

# Locate the Subtitle paragraph (holds "Comparative perspectives from Brazil, the UK, and the US").
$subtitle = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Subtitle") {
        $subtitle = $p
        break
    }
}

# Insert a new paragraph after the subtitle for the first author.
$subtitle.Range.InsertParagraphAfter()
$author1 = $subtitle.Next()
$author1.Range.Text = "Ben Jarman"
$author1.Style = "Author"

# Insert a new paragraph after the first author for the second author.
$author1.Range.InsertParagraphAfter()
$author2 = $author1.Next()
$author2.Range.Text = "Catherine Heard"
$author2.Style = "Author"

Write-Host "Inserted author paragraphs"
